$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("begroting")

# Remove the two obsolete sheets ("api calls" and "Sheet3") - the FORTRAN
# subroutine / API-call reference tables they held are no longer needed.
# Excel automatically drops the now-unused shared strings from the
# workbook's shared string table when the file is saved.
$wb.Worksheets("api calls").Delete() | Out-Null
$wb.Worksheets("Sheet3").Delete() | Out-Null

# "begroting" is now the only sheet - make sure it's the active/selected one.
$ws.Activate() | Out-Null

# Insert two blank rows above the "Total hours spent so far" summary row
# (currently row 36), pushing it down to row 38 and leaving room for two
# new line items (new test files / FindFacesRecursive / HangingEdge test).
$ws.Rows("36:37").Insert() | Out-Null

# Restore the on-screen selection to C30, matching the edited workbook.
$ws.Range("C30").Select() | Out-Null
